# Convert the boolean-typed FALSE cells in B3:G9 into plain numeric 0 cells.
# (Commit message: "new format for excel")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 9; $row++) {
    for ($col = 2; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        # Clear any boolean formatting/type by writing a pure numeric value.
        $cell.Value = 0
    }
}
